# Refresh Ultros market-board derived columns (currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ), LeveProfit(NQ/HQ)) for the leves whose Universalis
# price snapshot changed, one worksheet (Leve class) at a time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53: No Accounting for Waste
$ws.Cells.Item(53, 8).Value = 431.23077
$ws.Cells.Item(53, 9).Value = 179.5
$ws.Cells.Item(53, 11).Value = 179.5
$ws.Cells.Item(53, 13).Value = 457.5
# Row 88: The Grave of Hemlock Groves
$ws.Cells.Item(88, 8).Value = 4325.8184
$ws.Cells.Item(88, 10).Value = 1711.75
$ws.Cells.Item(88, 12).Value = 1711.75
$ws.Cells.Item(88, 14).Value = -2523.75
# Row 91: Dappling the Highlands (L)
$ws.Cells.Item(91, 8).Value = 4325.8184
$ws.Cells.Item(91, 10).Value = 1711.75
$ws.Cells.Item(91, 12).Value = 1711.75
$ws.Cells.Item(91, 14).Value = -4519.75
# Row 107: Another Man's Ink
$ws.Cells.Item(107, 8).Value = 1496.6111
$ws.Cells.Item(107, 9).Value = 1819.8462
$ws.Cells.Item(107, 10).Value = 656.2
$ws.Cells.Item(107, 11).Value = 1819.8462
$ws.Cells.Item(107, 12).Value = 656.2
$ws.Cells.Item(107, 13).Value = 100.1538
$ws.Cells.Item(107, 14).Value = -4496.2
# Row 134: Binding Spells
$ws.Cells.Item(134, 8).Value = 70000
$ws.Cells.Item(134, 10).Value = 70000
$ws.Cells.Item(134, 12).Value = 70000
$ws.Cells.Item(134, 14).Value = -80140
# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 2400.3845
$ws.Cells.Item(137, 10).Value = 2977.2222
$ws.Cells.Item(137, 12).Value = 8931.6666
$ws.Cells.Item(137, 14).Value = -14031.6666
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 5042.3696
$ws.Cells.Item(32, 9).Value = 5159.1855
$ws.Cells.Item(32, 11).Value = 5159.1855
$ws.Cells.Item(32, 13).Value = -4872.1855
# Row 74: As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 2291.2222
$ws.Cells.Item(74, 9).Value = 2376.1333
$ws.Cells.Item(74, 10).Value = 1866.6666
$ws.Cells.Item(74, 11).Value = 2376.1333
$ws.Cells.Item(74, 12).Value = 1866.6666
$ws.Cells.Item(74, 13).Value = -1502.1333
$ws.Cells.Item(74, 14).Value = -3614.6666
# Row 77: Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 2291.2222
$ws.Cells.Item(77, 9).Value = 2376.1333
$ws.Cells.Item(77, 10).Value = 1866.6666
$ws.Cells.Item(77, 11).Value = 11880.6665
$ws.Cells.Item(77, 12).Value = 9333.333000000001
$ws.Cells.Item(77, 13).Value = -7512.666499999999
$ws.Cells.Item(77, 14).Value = -18069.333
# Row 102: Smells of Rich Tama-hagane
$ws.Cells.Item(102, 8).Value = 4003.8667
$ws.Cells.Item(102, 9).Value = 4003.8667
$ws.Cells.Item(102, 11).Value = 4003.8667
$ws.Cells.Item(102, 13).Value = -2381.8667
# Row 105: Spoony Is the Bard
$ws.Cells.Item(105, 8).Value = 30000
$ws.Cells.Item(105, 10).Value = 30000
$ws.Cells.Item(105, 12).Value = 30000
$ws.Cells.Item(105, 14).Value = -36988
# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 1655.2354
$ws.Cells.Item(132, 9).Value = 1688.5
$ws.Cells.Item(132, 11).Value = 5065.5
$ws.Cells.Item(132, 13).Value = -2535.5
$ws = $wb.Worksheets.Item("BSM")
# Row 75: I Saw the Pine
$ws.Cells.Item(75, 8).Value = 130250
$ws.Cells.Item(75, 10).Value = 130250
$ws.Cells.Item(75, 12).Value = 130250
$ws.Cells.Item(75, 14).Value = -132122
# Row 78: I Came, I Sawed, I Conquered (L)
$ws.Cells.Item(78, 8).Value = 130250
$ws.Cells.Item(78, 10).Value = 130250
$ws.Cells.Item(78, 12).Value = 390750
$ws.Cells.Item(78, 14).Value = -400110
# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 3474.6
$ws.Cells.Item(94, 9).Value = 2731.1904
$ws.Cells.Item(94, 11).Value = 2731.1904
$ws.Cells.Item(94, 13).Value = -2280.1904
# Row 107: The Gold Experience
$ws.Cells.Item(107, 8).Value = 3512.1943
$ws.Cells.Item(107, 9).Value = 3109.8064
$ws.Cells.Item(107, 10).Value = 6007
$ws.Cells.Item(107, 11).Value = 3109.8064
$ws.Cells.Item(107, 12).Value = 6007
$ws.Cells.Item(107, 13).Value = -1189.8064
$ws.Cells.Item(107, 14).Value = -9847
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 1542.0588
$ws.Cells.Item(31, 9).Value = 1388.037
$ws.Cells.Item(31, 11).Value = 1388.037
$ws.Cells.Item(31, 13).Value = -1093.037
# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 1542.0588
$ws.Cells.Item(34, 9).Value = 1388.037
$ws.Cells.Item(34, 11).Value = 1388.037
$ws.Cells.Item(34, 13).Value = -1186.037
# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 1293.909
$ws.Cells.Item(58, 10).Value = 4700
$ws.Cells.Item(58, 12).Value = 4700
$ws.Cells.Item(58, 14).Value = -5106
# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 2416.5454
$ws.Cells.Item(122, 9).Value = 2664.652
$ws.Cells.Item(122, 10).Value = 1845.9
$ws.Cells.Item(122, 11).Value = 7993.956
$ws.Cells.Item(122, 12).Value = 5537.700000000001
$ws.Cells.Item(122, 13).Value = -5543.956
$ws.Cells.Item(122, 14).Value = -10437.7
# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 2697.3215
$ws.Cells.Item(132, 9).Value = 2008.1111
$ws.Cells.Item(132, 10).Value = 3937.9
$ws.Cells.Item(132, 11).Value = 6024.3333
$ws.Cells.Item(132, 12).Value = 11813.7
$ws.Cells.Item(132, 13).Value = -3494.3333
$ws.Cells.Item(132, 14).Value = -16873.7
# Row 134: Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 3379.611
$ws.Cells.Item(134, 9).Value = 3434.6667
$ws.Cells.Item(134, 10).Value = 3104.3333
$ws.Cells.Item(134, 11).Value = 10304.0001
$ws.Cells.Item(134, 12).Value = 9312.999899999999
$ws.Cells.Item(134, 13).Value = -7769.000100000001
$ws.Cells.Item(134, 14).Value = -14382.9999
# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 1293.909
$ws.Cells.Item(136, 10).Value = 4700
$ws.Cells.Item(136, 12).Value = 14100
$ws.Cells.Item(136, 14).Value = -19200
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Cells.Item(5, 8).Value = 944.85297
$ws.Cells.Item(5, 9).Value = 977.9643
$ws.Cells.Item(5, 11).Value = 2933.8929
$ws.Cells.Item(5, 13).Value = -2821.8929
# Row 37: I Love Lamprey
$ws.Cells.Item(37, 8).Value = 649859.75
$ws.Cells.Item(37, 10).Value = 649859.75
$ws.Cells.Item(37, 12).Value = 1949579.25
$ws.Cells.Item(37, 14).Value = -1949803.25
# Row 106: Herky Jerky
$ws.Cells.Item(106, 8).Value = 5000
$ws.Cells.Item(106, 10).Value = 5000
$ws.Cells.Item(106, 12).Value = 15000
$ws.Cells.Item(106, 14).Value = -16892
# Row 135: Not-so-secret Ingredient
$ws.Cells.Item(135, 8).Value = 944.85297
$ws.Cells.Item(135, 9).Value = 977.9643
$ws.Cells.Item(135, 11).Value = 8801.6787
$ws.Cells.Item(135, 13).Value = -6266.6787
# Row 136: Simple Is Hardest
$ws.Cells.Item(136, 8).Value = 3814.2222
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 14).Value = -25200
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 43592.2
$ws.Cells.Item(80, 9).Value = 66705.11
$ws.Cells.Item(80, 10).Value = 8922.833000000001
$ws.Cells.Item(80, 11).Value = 66705.11
$ws.Cells.Item(80, 12).Value = 8922.833000000001
$ws.Cells.Item(80, 13).Value = -65707.11
$ws.Cells.Item(80, 14).Value = -10918.833
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 43592.2
$ws.Cells.Item(83, 9).Value = 66705.11
$ws.Cells.Item(83, 10).Value = 8922.833000000001
$ws.Cells.Item(83, 11).Value = 333525.55
$ws.Cells.Item(83, 12).Value = 44614.165
$ws.Cells.Item(83, 13).Value = -328533.55
$ws.Cells.Item(83, 14).Value = -54598.165
# Row 113: Copious Crystal Cannons
$ws.Cells.Item(113, 8).Value = 5316.8887
$ws.Cells.Item(113, 9).Value = 1121.7142
$ws.Cells.Item(113, 11).Value = 1121.7142
$ws.Cells.Item(113, 13).Value = 1048.2858
$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck
$ws.Cells.Item(82, 8).Value = 66668190
$ws.Cells.Item(82, 10).Value = 1149.75
$ws.Cells.Item(82, 12).Value = 1149.75
$ws.Cells.Item(82, 14).Value = -1871.75
# Row 85: Training Is Only Skintight (L)
$ws.Cells.Item(85, 8).Value = 66668190
$ws.Cells.Item(85, 10).Value = 1149.75
$ws.Cells.Item(85, 12).Value = 1149.75
$ws.Cells.Item(85, 14).Value = -3645.75
# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 2385.1143
$ws.Cells.Item(136, 9).Value = 2181.6897
$ws.Cells.Item(136, 10).Value = 3368.3333
$ws.Cells.Item(136, 11).Value = 6545.0691
$ws.Cells.Item(136, 12).Value = 10104.9999
$ws.Cells.Item(136, 13).Value = -3995.0691
$ws.Cells.Item(136, 14).Value = -15204.9999
$ws = $wb.Worksheets.Item("WVR")
# Row 45: Private Concerns
$ws.Cells.Item(45, 8).Value = 14651.625
$ws.Cells.Item(45, 9).Value = 7783
$ws.Cells.Item(45, 10).Value = 16941.166
$ws.Cells.Item(45, 11).Value = 7783
$ws.Cells.Item(45, 12).Value = 16941.166
$ws.Cells.Item(45, 13).Value = -7292
$ws.Cells.Item(45, 14).Value = -17923.166
# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 2130.5454
$ws.Cells.Item(132, 10).Value = 5500
$ws.Cells.Item(132, 12).Value = 16500
$ws.Cells.Item(132, 14).Value = -21560
